$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price report row is published; it becomes the new row 2 and
# every existing data row (formerly 2-10) shifts down by one (now 3-11).
$ws.Rows.Item(2).Insert()

# The freshly inserted row inherits the header row's bold/centered style from
# the insert; strip that back to the default (unstyled) look used by the
# other data rows before we fill it in.
$ws.Range("A2:T2").ClearFormats()

$ws.Cells.Item(2, 1).Value = 11
$ws.Cells.Item(2, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(2, 3).Value = "Bíobío"
$ws.Cells.Item(2, 4).Value = 44425
$ws.Cells.Item(2, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(2, 5).Value = 8
$ws.Cells.Item(2, 6).Value = "Fruta"
$ws.Cells.Item(2, 7).Value = 100104
$ws.Cells.Item(2, 8).Value = "Frutos de pepita"
$ws.Cells.Item(2, 9).Value = 100104003
$ws.Cells.Item(2, 10).Value = "Membrillo"
$ws.Cells.Item(2, 11).Value = "Champion"
$ws.Cells.Item(2, 12).Value = "Primera"
$ws.Cells.Item(2, 13).Value = 100
$ws.Cells.Item(2, 14).Value = 12000
$ws.Cells.Item(2, 15).Value = 13000
$ws.Cells.Item(2, 16).Value = 12500
$ws.Cells.Item(2, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(2, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(2, 19).Value = 694
$ws.Cells.Item(2, 20).Value = 18
